$wb = $excel.ActiveWorkbook

# 1. Status text change: "Ready for handoff" -> "In Translation"
#    This string appears on the Overview sheet (zh-cn / de-de status columns)
#    and on each per-locale sheet's "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2. Narrow the "Status" columns to match the shorter text now shown
#    (was sized for "Ready for handoff", now resized for "In Translation").
$wsOverview.Columns.Item(5).ColumnWidth = 12.58
$wsOverview.Columns.Item(6).ColumnWidth = 12.58

$wsZhCn.Columns.Item(3).ColumnWidth = 12.58

$wsDeDe.Columns.Item(3).ColumnWidth = 12.58
